# Fruta / hortaliza, semanal
# Rotates the per-row data (Fecha, Volumen, Precio mínimo/máximo/promedio,
# Unidad de comercialización, Precio $/Kg) among rows 2-9, leaving all other
# columns (which are identical down the sheet) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each target row, the tuple of values that should end up there (taken
# from another row in the original layout).
$data = @{
    2 = @(44208, 210, 10000, 10000, 10000, "`$/caja 14 kilos empedrada", 714)
    3 = @(44176, 250, 7000,  7000,  7000,  "`$/caja 14 kilos empedrada", 500)
    4 = @(44351, 300, 10000, 10000, 10000, "`$/caja 14 kilos empedrada", 714)
    5 = @(44162, 120, 7000,  7000,  7000,  "`$/caja 14 kilos empedrada", 500)
    6 = @(44491, 180, 9000,  9000,  9000,  "`$/caja 14 kilos empedrada", 643)
    7 = @(44397, 60,  11000, 11000, 11000, "`$/caja 14 kilos",           786)
    8 = @(44400, 100, 10000, 10000, 10000, "`$/caja 14 kilos",           714)
    9 = @(44309, 300, 7000,  7000,  7000,  "`$/caja 14 kilos empedrada", 500)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]

    $ws.Range("D$row").Value = $values[0]
    $ws.Range("M$row").Value = $values[1]
    $ws.Range("N$row").Value = $values[2]
    $ws.Range("O$row").Value = $values[3]
    $ws.Range("P$row").Value = $values[4]
    $ws.Range("Q$row").Value = $values[5]
    $ws.Range("S$row").Value = $values[6]
}
